$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.010.62'
$ws.Range('E2').Value = '  -0.03%  '
$ws.Range('D3').Value = '2.234.89'
$ws.Range('E3').Value = '  -0.45%  '
$ws.Range('E4').Value = '  +0.14%  '
$st = $ws.Range('D5').Style
$ws.Range('D5').Value = "'304.54"
$ws.Range('D5').Style = $st
$ws.Range('E5').Value = '  -4.31%  '
$st = $ws.Range('D6').Style
$ws.Range('D6').Value = "'93.93"
$ws.Range('D6').Style = $st
$ws.Range('E6').Value = '  -6.80%  '
$st = $ws.Range('D7').Style
$ws.Range('D7').Value = "'0.569"
$ws.Range('D7').Style = $st
$ws.Range('E7').Value = '  -0.57%  '
$ws.Range('E8').Value = '  +0.19%  '
$st = $ws.Range('D9').Style
$ws.Range('D9').Value = "'0.521"
$ws.Range('D9').Style = $st
$ws.Range('E9').Value = '  -3.94%  '
$st = $ws.Range('D10').Style
$ws.Range('D10').Value = "'34.64"
$ws.Range('D10').Style = $st
$ws.Range('E10').Value = '  -6.04%  '
$ws.Range('E11').Value = '  -2.18%  '
$ws.Range('E12').Value = '  -4.47%  '
$ws.Range('E13').Value = '  -0.80%  '
$ws.Range('D14').Value = '2.575.80'
$ws.Range('E14').Value = '  -0.50%  '
$ws.Range('D15').Value = '2.236.35'
$ws.Range('E15').Value = '  -0.43%  '
$st = $ws.Range('D16').Style
$ws.Range('D16').Value = "'0.819"
$ws.Range('D16').Style = $st
$ws.Range('E16').Value = '  -3.14%  '
$st = $ws.Range('D17').Style
$ws.Range('D17').Value = "'13.47"
$ws.Range('D17').Style = $st
$ws.Range('E17').Value = '  -4.87%  '
$ws.Range('D18').Value = '43.870.69'
$ws.Range('E18').Value = '  -0.20%  '
$ws.Range('D19').Value = '0.0₃0960'
$ws.Range('E19').Value = '  -1.48%  '
$st = $ws.Range('D20').Style
$ws.Range('D20').Value = "'12.06"
$ws.Range('D20').Style = $st
$ws.Range('E20').Value = '  -9.39%  '
$ws.Range('E21').Value = '  -2.16%  '
$st = $ws.Range('D22').Style
$ws.Range('D22').Value = "'65.44"
$ws.Range('D22').Style = $st
$ws.Range('E22').Value = '  -0.02%  '
$st = $ws.Range('D23').Style
$ws.Range('D23').Value = "'236.38"
$ws.Range('D23').Style = $st
$ws.Range('E23').Value = '  +0.73%  '
$st = $ws.Range('D24').Style
$ws.Range('D24').Value = "'2.90"
$ws.Range('D24').Style = $st
$ws.Range('E24').Value = '  -6.24%  '
$ws.Range('E25').Value = '  -4.90%  '
$ws.Range('E26').Value = '  +0.16%  '
$ws.Range('B27').Value = 'Toncoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$st = $ws.Range('D27').Style
$ws.Range('D27').Value = "'2.20"
$ws.Range('D27').Style = $st
$ws.Range('E27').Value = '  -0.36%  '
$ws.Range('B28').Value = 'InjectiveProtocol'
$ws.Range('C28').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$st = $ws.Range('D28').Style
$ws.Range('D28').Value = "'37.98"
$ws.Range('D28').Style = $st
$ws.Range('E28').Value = '  +0.44%  '
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$st = $ws.Range('D29').Style
$ws.Range('D29').Value = "'9.84"
$ws.Range('D29').Style = $st
$ws.Range('E29').Value = '  -6.02%  '
$st = $ws.Range('D30').Style
$ws.Range('D30').Value = "'6.01"
$ws.Range('D30').Style = $st
$ws.Range('E30').Value = '  -1.66%  '
$st = $ws.Range('D31').Style
$ws.Range('D31').Value = "'19.84"
$ws.Range('D31').Style = $st
$ws.Range('E31').Value = '  -1.05%  '
$st = $ws.Range('D32').Style
$ws.Range('D32').Value = "'150.30"
$ws.Range('D32').Style = $st
$ws.Range('E32').Value = '  -5.65%  '
$st = $ws.Range('D33').Style
$ws.Range('D33').Value = "'0.0795"
$ws.Range('D33').Style = $st
$ws.Range('E33').Value = '  -5.78%  '
$ws.Range('E34').Value = '  -3.06%  '
$st = $ws.Range('D35').Style
$ws.Range('D35').Value = "'3.15"
$ws.Range('D35').Style = $st
$ws.Range('E35').Value = '  -1.61%  '
$ws.Range('E36').Value = '  -3.54%  '
$ws.Range('E37').Value = '  +0.89%  '
$ws.Range('E38').Value = '  -8.87%  '
$st = $ws.Range('D39').Style
$ws.Range('D39').Value = "'14.88"
$ws.Range('D39').Style = $st
$ws.Range('E39').Value = '  -7.14%  '
$st = $ws.Range('D40').Style
$ws.Range('D40').Value = "'3.83"
$ws.Range('D40').Style = $st
$ws.Range('E40').Value = '  -7.50%  '
$st = $ws.Range('D41').Style
$ws.Range('D41').Value = "'3.36"
$ws.Range('D41').Style = $st
$ws.Range('E41').Value = '  -8.62%  '
$ws.Range('E42').Value = '  -6.05%  '
$ws.Range('E43').Value = '  +0.16%  '
$ws.Range('D44').Value = '1.732.22'
$ws.Range('E44').Value = '  -0.57%  '
$st = $ws.Range('D45').Style
$ws.Range('D45').Value = "'84.56"
$ws.Range('D45').Style = $st
$ws.Range('E45').Value = '  +3.58%  '
$st = $ws.Range('D46').Style
$ws.Range('D46').Value = "'0.186"
$ws.Range('D46').Style = $st
$ws.Range('E46').Value = '  -5.29%  '
$st = $ws.Range('D47').Style
$ws.Range('D47').Value = "'99.54"
$ws.Range('D47').Style = $st
$ws.Range('E47').Value = '  -2.64%  '
$ws.Range('E48').Value = '  -4.61%  '
$st = $ws.Range('D49').Style
$ws.Range('D49').Value = "'8.05"
$ws.Range('D49').Style = $st
$ws.Range('E49').Value = '  -2.14%  '
$st = $ws.Range('D50').Style
$ws.Range('D50').Value = "'68.44"
$ws.Range('D50').Style = $st
$ws.Range('E50').Value = '  -7.91%  '
$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$st = $ws.Range('D51').Style
$ws.Range('D51').Value = "'53.69"
$ws.Range('D51').Style = $st
$ws.Range('E51').Value = '  -6.54%  '
